$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.951.86'
$ws.Range('E2').Value = '  +6.51%  '
$ws.Range('D3').Value = '2.643.75'
$ws.Range('E3').Value = '  +8.85%  '
$ws.Range('D4').Value = "'1.01"
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').Value = "'514.63"
$ws.Range('E5').Value = '  +4.91%  '
$ws.Range('D6').Value = "'159.69"
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('E7').Value = '  -0.86%  '
$ws.Range('D8').Value = "'0.997"
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = '2.685.94'
$ws.Range('E9').Value = '  +9.62%  '
$ws.Range('D10').Value = "'6.18"
$ws.Range('E10').Value = '  +8.77%  '
$ws.Range('E11').Value = '  +5.79%  '
$ws.Range('E12').Value = '  +3.76%  '
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').Value = '3.119.34'
$ws.Range('E14').Value = '  +9.37%  '
$ws.Range('D15').Value = '61.038.98'
$ws.Range('E15').Value = '  +6.66%  '
$ws.Range('D16').Value = "'22.30"
$ws.Range('E16').Value = '  +6.73%  '
$ws.Range('E17').Value = '  +5.01%  '
$ws.Range('D18').Value = '2.686.15'
$ws.Range('E18').Value = '  +9.74%  '
$ws.Range('D19').Value = "'4.84"
$ws.Range('E19').Value = '  +1.38%  '
$ws.Range('D20').Value = "'349.92"
$ws.Range('E20').Value = '  +6.09%  '
$ws.Range('D21').Value = "'10.57"
$ws.Range('E21').Value = '  +5.70%  '
$ws.Range('E22').Value = '  +4.59%  '
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').Value = "'60.77"
$ws.Range('E24').Value = '  +4.25%  '
$ws.Range('E25').Value = '  +3.46%  '
$ws.Range('D26').Value = '2.756.28'
$ws.Range('E26').Value = '  +8.54%  '
$ws.Range('E27').Value = '  +4.71%  '
$ws.Range('D28').Value = "'0.994"
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').Value = '0.0₃0873'
$ws.Range('E29').Value = '  +10.00%  '
$ws.Range('D30').Value = "'7.59"
$ws.Range('E30').Value = '  +3.62%  '
$ws.Range('D31').Value = "'1.00"
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').Value = "'19.72"
$ws.Range('E32').Value = '  +4.83%  '
$ws.Range('D33').Value = "'157.66"
$ws.Range('E33').Value = '  +5.56%  '
$ws.Range('E34').Value = '  +4.75%  '
$ws.Range('E35').Value = '  +7.45%  '
$ws.Range('E36').Value = '  +10.16%  '
$ws.Range('E37').Value = '  +5.97%  '
$ws.Range('D38').Value = "'0.888"
$ws.Range('E38').Value = '  +3.30%  '
$ws.Range('D39').Value = "'1.55"
$ws.Range('E39').Value = '  +11.63%  '
$ws.Range('D40').Value = "'311.96"
$ws.Range('E40').Value = '  +16.51%  '
$ws.Range('D41').Value = "'3.81"
$ws.Range('E41').Value = '  +7.55%  '
$ws.Range('D42').Value = "'0.839"
$ws.Range('E42').Value = '  +30.06%  '
$ws.Range('D43').Value = "'35.74"
$ws.Range('E43').Value = '  +4.34%  '
$ws.Range('D44').Value = "'0.650"
$ws.Range('E44').Value = '  +8.59%  '
$ws.Range('D45').Value = "'0.0580"
$ws.Range('E45').Value = '  +7.92%  '
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').Value = "'20.36"
$ws.Range('E47').Value = '  +15.76%  '
$ws.Range('D48').Value = "'0.991"
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('D49').Value = "'5.02"
$ws.Range('E49').Value = '  +6.83%  '
$ws.Range('D50').Value = "'0.0238"
$ws.Range('E50').Value = '  +3.48%  '
$ws.Range('D51').Value = '2.042.03'
$ws.Range('E51').Value = '  +9.63%  '

$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
